$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 6
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 7
$ws.Range("B8").Value = "<number>"
$ws.Range("C9").Value = 2
$ws.Range("C10").Value = 4
$ws.Range("C12").Value = 5
$ws.Range("B13").Value = "<tho>"
$ws.Range("C13").Value = 8
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 4
